$d = $word.ActiveDocument

$replacements = @(
    @("786×3=2358", "486×7=3402"),
    @("610×3=1830", "520×5=2600"),
    @("786×4=3144", "667×6=4002"),
    @("830×9=7470", "953×6=5718"),
    @("224×5=1120", "225×5=1125"),
    @("465×7=3255", "965×9=8685"),
    @("856×4=3424", "920×6=5520"),
    @("791×2=1582", "684×6=4104"),
    @("256×5=1280", "310×9=2790"),
    @("300×4=1200", "972×7=6804"),
    @("722×4=2888", "852×8=6816"),
    @("607×8=4856", "386×4=1544"),
    @("910×7=6370", "245×7=1715"),
    @("348×5=1740", "894×6=5364"),
    @("729×5=3645", "139×4=556"),
    @("590×8=4720", "623×9=5607"),
    @("253×5=1265", "121×5=605"),
    @("409×8=3272", "177×8=1416"),
    @("227×2=454", "690×9=6210"),
    @("735×9=6615", "373×8=2984"),
    @("945×2=1890", "975×2=1950"),
    @("632×8=5056", "810×6=4860"),
    @("603×3=1809", "345×6=2070"),
    @("552×2=1104", "864×8=6912"),
    @("737×9=6633", "557×5=2785")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
